$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 1718.98
$ws.Range("I15").Value2 = 1718.98
$ws.Range("K15").Value2 = 5156.940000000001
$ws.Range("M15").Value2 = -4987.940000000001
$ws.Range("H31").Value2 = 11357
$ws.Range("I31").Value2 = 11357
$ws.Range("K31").Value2 = 34071
$ws.Range("M31").Value2 = -33841
$ws.Range("H39").Value2 = 280.72
$ws.Range("I39").Value2 = 98.9375
$ws.Range("J39").Value2 = 603.8889
$ws.Range("K39").Value2 = 296.8125
$ws.Range("L39").Value2 = 1811.6667
$ws.Range("M39").Value2 = -0.8125
$ws.Range("N39").Value2 = -2403.6667
$ws.Range("H98").Value2 = 7514.067
$ws.Range("I98").Value2 = 4693.154
$ws.Range("K98").Value2 = 4693.154
$ws.Range("M98").Value2 = -3195.154
$ws.Range("H103").Value2 = 426.25
$ws.Range("I103").Value2 = 353.33334
$ws.Range("J103").Value2 = 645
$ws.Range("K103").Value2 = 1060.00002
$ws.Range("L103").Value2 = 1935
$ws.Range("M103").Value2 = -474.0000199999999
$ws.Range("N103").Value2 = -3107
$ws.Range("H106").Value2 = 6300
$ws.Range("I106").Value2 = 5950
$ws.Range("K106").Value2 = 5950
$ws.Range("M106").Value2 = -5319
$ws.Range("H107").Value2 = 565.125
$ws.Range("I107").Value2 = 467.5
$ws.Range("J107").Value2 = 662.75
$ws.Range("K107").Value2 = 467.5
$ws.Range("L107").Value2 = 662.75
$ws.Range("M107").Value2 = 1452.5
$ws.Range("N107").Value2 = -4502.75
$ws.Range("H112").Value2 = 1988.619
$ws.Range("I112").Value2 = 900
$ws.Range("J112").Value2 = 2103.2104
$ws.Range("K112").Value2 = 2700
$ws.Range("L112").Value2 = 6309.6312
$ws.Range("M112").Value2 = -1592
$ws.Range("N112").Value2 = -8525.6312
$ws.Range("H122").Value2 = 7514.067
$ws.Range("I122").Value2 = 4693.154
$ws.Range("K122").Value2 = 14079.462
$ws.Range("M122").Value2 = -11629.462
$ws.Range("H127").Value2 = 797.3333
$ws.Range("I127").Value2 = 797.3333
$ws.Range("J127").Value2 = 0
$ws.Range("K127").Value2 = 2391.9999
$ws.Range("L127").Value2 = 0
$ws.Range("M127").Value2 = 2568.0001
$ws.Range("N127").Value2 = $null

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 1797.3334
$ws.Range("I74").Value2 = 3375
$ws.Range("J74").Value2 = 1008.5
$ws.Range("K74").Value2 = 3375
$ws.Range("L74").Value2 = 1008.5
$ws.Range("M74").Value2 = -2501
$ws.Range("N74").Value2 = -2756.5
$ws.Range("H77").Value2 = 1797.3334
$ws.Range("I77").Value2 = 3375
$ws.Range("J77").Value2 = 1008.5
$ws.Range("K77").Value2 = 16875
$ws.Range("L77").Value2 = 5042.5
$ws.Range("M77").Value2 = -12507
$ws.Range("N77").Value2 = -13778.5
$ws.Range("H123").Value2 = 23085
$ws.Range("J123").Value2 = 23085
$ws.Range("L123").Value2 = 23085
$ws.Range("N123").Value2 = -32885

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 44131.707
$ws.Range("I86").Value2 = 2556.9473
$ws.Range("J86").Value2 = 202115.8
$ws.Range("K86").Value2 = 2556.9473
$ws.Range("L86").Value2 = 202115.8
$ws.Range("M86").Value2 = -1433.9473
$ws.Range("N86").Value2 = -204361.8
$ws.Range("H89").Value2 = 44131.707
$ws.Range("I89").Value2 = 2556.9473
$ws.Range("J89").Value2 = 202115.8
$ws.Range("K89").Value2 = 12784.7365
$ws.Range("L89").Value2 = 1010579
$ws.Range("M89").Value2 = -7168.736499999999
$ws.Range("N89").Value2 = -1021811

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 416.52942
$ws.Range("I22").Value2 = 241.07692
$ws.Range("J22").Value2 = 986.75
$ws.Range("K22").Value2 = 241.07692
$ws.Range("L22").Value2 = 986.75
$ws.Range("M22").Value2 = 108.92308
$ws.Range("N22").Value2 = -1686.75
$ws.Range("H88").Value2 = 38000
$ws.Range("J88").Value2 = 38000
$ws.Range("L88").Value2 = 38000
$ws.Range("N88").Value2 = -38812
$ws.Range("H91").Value2 = 38000
$ws.Range("J91").Value2 = 38000
$ws.Range("L91").Value2 = 38000
$ws.Range("N91").Value2 = -40808
$ws.Range("H122").Value2 = 2349.75
$ws.Range("I122").Value2 = 1799.6666
$ws.Range("J122").Value2 = 4000
$ws.Range("K122").Value2 = 5398.9998
$ws.Range("L122").Value2 = 12000
$ws.Range("M122").Value2 = -2948.9998
$ws.Range("N122").Value2 = -16900

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value2 = 192950.64
$ws.Range("I113").Value2 = 294691.6
$ws.Range("J113").Value2 = 773.2778
$ws.Range("K113").Value2 = 884074.7999999999
$ws.Range("L113").Value2 = 2319.8334
$ws.Range("M113").Value2 = -881904.7999999999
$ws.Range("N113").Value2 = -6659.8334
$ws.Range("H136").Value2 = 6137.8667
$ws.Range("J136").Value2 = 6362
$ws.Range("L136").Value2 = 19086
$ws.Range("N136").Value2 = -29286

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value2 = 6707.2856
$ws.Range("J92").Value2 = 6707.2856
$ws.Range("L92").Value2 = 6707.2856
$ws.Range("N92").Value2 = -10451.2856
$ws.Range("H109").Value2 = 20284
$ws.Range("J109").Value2 = 20284
$ws.Range("L109").Value2 = 20284
$ws.Range("N109").Value2 = -22364
$ws.Range("H123").Value2 = 8938.866
$ws.Range("J123").Value2 = 8938.866
$ws.Range("L123").Value2 = 8938.866
$ws.Range("N123").Value2 = -13838.866
$ws.Range("H132").Value2 = 3429.8064
$ws.Range("I132").Value2 = 2231.2
$ws.Range("J132").Value2 = 4000.5715
$ws.Range("K132").Value2 = 6693.599999999999
$ws.Range("L132").Value2 = 12001.7145
$ws.Range("M132").Value2 = -4163.599999999999
$ws.Range("N132").Value2 = -17061.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 40002780
$ws.Range("I122").Value2 = 3004
$ws.Range("J122").Value2 = 50002724
$ws.Range("K122").Value2 = 9012
$ws.Range("L122").Value2 = 150008172
$ws.Range("M122").Value2 = -6562
$ws.Range("N122").Value2 = -150013072
$ws.Range("H132").Value2 = 4541.107
$ws.Range("I132").Value2 = 4564.1113
$ws.Range("J132").Value2 = 4499.7
$ws.Range("K132").Value2 = 13692.3339
$ws.Range("L132").Value2 = 13499.1
$ws.Range("M132").Value2 = -11162.3339
$ws.Range("N132").Value2 = -18559.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value2 = 29995
$ws.Range("J82").Value2 = 29995
$ws.Range("L82").Value2 = 29995
$ws.Range("N82").Value2 = -30761
$ws.Range("H85").Value2 = 29995
$ws.Range("J85").Value2 = 29995
$ws.Range("L85").Value2 = 29995
$ws.Range("N85").Value2 = -32647
$ws.Range("H123").Value2 = 21425.25
$ws.Range("J123").Value2 = 21425.25
$ws.Range("L123").Value2 = 21425.25
$ws.Range("N123").Value2 = -31225.25
$ws.Range("H125").Value2 = 60602.5
$ws.Range("J125").Value2 = 60602.5
$ws.Range("L125").Value2 = 60602.5
$ws.Range("N125").Value2 = -70442.5
$ws.Range("H126").Value2 = 7826.9165
$ws.Range("J126").Value2 = 4169.5713
$ws.Range("L126").Value2 = 12508.7139
$ws.Range("N126").Value2 = -17448.7139
